$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("date", "timestamp in format YYYY-MM-DD"),
    @("CountryCode", "standart ISO country code nomenclature of 3 characters"),
    @("StringencyIndex", "geverment response indicator for the covid prevention. More info on how this value is calculated can be found in the metadat form Oxford University"),
    @("Country_Region", "Country name"),
    @("Population", "total population taken from the metadata on Covid data source repository. More info in the JH repository metadata location"),
    @("Confirmed", "confirmed cases of covid infections"),
    @("Deaths", "confirmed deaths cases caused by covid infections"),
    @("Recovered", "death cases normalized by 100000 people"),
    @("Active", "current active cases of covid infections. Calculated by subtracting Recovered and Death cases to Confirmed cases"),
    @("New_cases", "lag difference of confirmed cases in day-by-day-basis"),
    @("New_deaths", "lag difference of death cases in day-by-day-basis"),
    @("New_recovered", "lag difference of recovered cases in day-by-day-basis"),
    @("Confirmed_100K", "confirmed cases normalized by 100000 people"),
    @("Deaths_100K", "death cases normalized by 100000 people"),
    @("Recovered_100K", "recovered cases normalized by 100000 people"),
    @("Active_100K", "active cases normalized by 100000 people"),
    @("New_cases_100K", "new cases normalized by 100000 people"),
    @("New_deaths_100K", "new death normalized by 100000 people"),
    @("New_recovered_100K", "new recovered normalized by 100000 people"),
    @("GDP_in_USD", "Gross domestic product of the country in USD currency"),
    @("incomeLevel.value", "Income level assigned to the country from the world bank in categorical  ")
)

$r = 2
foreach ($pair in $data) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $r++
}

# Column B is populated in the same order the shared-strings table shows them
# being created: rows 2-8, then 10-14, then 9, then 15, then 16-22.
$bOrder = @(2,3,4,5,6,7,8,10,11,12,13,14,9,15,16,17,18,19,20,21,22)
foreach ($rowNum in $bOrder) {
    $pair = $data[$rowNum - 2]
    $ws.Cells.Item($rowNum, 2).Value = $pair[1]
}

$ws.Columns.Item(1).ColumnWidth = 18.5
$ws.Columns.Item(2).ColumnWidth = 71.5

$ws.Range("B23").Select()

$activeWindow = $excel.ActiveWindow
$activeWindow.Zoom = 231
